$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 14:05"

# Update country names for rows whose ranking changed
$ws.Range("A47").Value = "Dinamarca"
$ws.Range("A48").Value = "Argentina"
$ws.Range("A53").Value = "Barein"
$ws.Range("A54").Value = "Chequia"
$ws.Range("A78").Value = "Senegal"
$ws.Range("A79").Value = "Tailandia"
$ws.Range("A80").Value = "Tayikistan"
$ws.Range("A81").Value = "Grecia"
$ws.Range("A134").Value = "Madagascar"
$ws.Range("A135").Value = "Tanzania"
$ws.Range("A198").Value = "Santa Lucia"
$ws.Range("A199").Value = "Belice"
$ws.Range("A210").Value = "Montserrat"
$ws.Range("A211").Value = "Groenlandia"
$ws.Range("A214").Value = "San Bartolome"
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "Sahara Occidental"

# Update numeric stats that changed
$ws.Range("B4").Value = 1666923
$ws.Range("C4").Value = 95
$ws.Range("E4").Value = 1121313
$ws.Range("B14").Value = 132755
$ws.Range("C14").Value = 1332
$ws.Range("D14").Value = 54819
$ws.Range("E14").Value = 74037
$ws.Range("B23").Value = 45236
$ws.Range("C23").Value = 172
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 5822
$ws.Range("B28").Value = 33459
$ws.Range("C28").Value = 271
$ws.Range("E28").Value = 24490
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 3998
$ws.Range("E30").Value = 830
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 1906
$ws.Range("B47").Value = 11360
$ws.Range("C47").Value = 71
$ws.Range("D47").Value = 9900
$ws.Range("E47").Value = 898
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 562
$ws.Range("B48").Value = 11353
$ws.Range("D48").Value = 3530
$ws.Range("E48").Value = 7378
$ws.Range("H48").Value = 445
$ws.Range("B53").Value = 9093
$ws.Range("C53").Value = 291
$ws.Range("D53").Value = 4581
$ws.Range("E53").Value = 4499
$ws.Range("H53").Value = 13
$ws.Range("B54").Value = 8891
$ws.Range("C54").Value = 1
$ws.Range("D54").Value = 6047
$ws.Range("E54").Value = 2529
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 315
$ws.Range("B76").Value = 3150
$ws.Range("C76").Value = 35
$ws.Range("D76").Value = 2565
$ws.Range("E76").Value = 572
$ws.Range("B78").Value = 3047
$ws.Range("C78").Value = 71
$ws.Range("D78").Value = 1456
$ws.Range("E78").Value = 1556
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 35
$ws.Range("B79").Value = 3040
$ws.Range("D79").Value = 2921
$ws.Range("E79").Value = 63
$ws.Range("H79").Value = 56
$ws.Range("B80").Value = 2929
$ws.Range("C80").Value = 191
$ws.Range("D80").Value = 1301
$ws.Range("E80").Value = 1582
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 46
$ws.Range("B81").Value = 2876
$ws.Range("D81").Value = 1374
$ws.Range("E81").Value = 1331
$ws.Range("H81").Value = 171
$ws.Range("B134").Value = 527
$ws.Range("C134").Value = 39
$ws.Range("D134").Value = 138
$ws.Range("E134").Value = 387
$ws.Range("H134").Value = 2
$ws.Range("B135").Value = 509
$ws.Range("D135").Value = 183
$ws.Range("E135").Value = 305
$ws.Range("H135").Value = 21
$ws.Range("D139").Value = 357
$ws.Range("E139").Value = 63
$ws.Range("D198").Value = 18
$ws.Range("H198").Value = 0
$ws.Range("D199").Value = 16
$ws.Range("H199").Value = 2
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0
